# Update the fixed "datetimeFigureOut" date placeholder text from
# 12/25/2020 to 12/27/2020 everywhere it appears: the slide master and
# every slide layout (the actual slide(s) don't carry their own copy of
# this placeholder text).

$p = $ppt.ActivePresentation

$oldDate = "12/25/2020"
$newDate = "12/27/2020"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -eq -1 -and $shp.Type -eq 14) {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# Slide master.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout hanging off the master.
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholder $layouts.Item($L).Shapes
}
